$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PDD")

$ws.Range("B6").Value = 6.1982
$ws.Range("B7").Value = 6.2567
$ws.Range("B8").Value = 3.4618
$ws.Range("B9").Value = 47.1575
$ws.Range("B10").Value = 11.2995
$ws.Range("B11").Value = 12.7896
$ws.Range("B38").Value = 14.3277
$ws.Range("B58").Value = 47.1575
